$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 33.57143
$ws.Range("I38").Value = 33.57143
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 100.71429
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 271.28571
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 7132.8
$ws.Range("I40").Value = 2998.6
$ws.Range("J40").Value = 9199.9
$ws.Range("K40").Value = 2998.6
$ws.Range("L40").Value = 9199.9
$ws.Range("M40").Value = -2823.6
$ws.Range("N40").Value = -9549.9
$ws.Range("H51").Value = 6379.4
$ws.Range("J51").Value = 1000
$ws.Range("L51").Value = 1000
$ws.Range("N51").Value = -1968
$ws.Range("H75").Value = 83437
$ws.Range("J75").Value = 83437
$ws.Range("L75").Value = 83437
$ws.Range("N75").Value = -85309
$ws.Range("H78").Value = 83437
$ws.Range("J78").Value = 83437
$ws.Range("L78").Value = 250311
$ws.Range("N78").Value = -259671
$ws.Range("H98").Value = 1737.5454
$ws.Range("I98").Value = 1841.8
$ws.Range("J98").Value = 695
$ws.Range("K98").Value = 1841.8
$ws.Range("L98").Value = 695
$ws.Range("M98").Value = -343.8
$ws.Range("N98").Value = -3691
$ws.Range("H100").Value = 2782.8333
$ws.Range("J100").Value = 1997.5
$ws.Range("L100").Value = 1997.5
$ws.Range("N100").Value = -3079.5
$ws.Range("H107").Value = 667.7059
$ws.Range("J107").Value = 1062.75
$ws.Range("L107").Value = 1062.75
$ws.Range("N107").Value = -4902.75
$ws.Range("H122").Value = 1737.5454
$ws.Range("I122").Value = 1841.8
$ws.Range("J122").Value = 695
$ws.Range("K122").Value = 5525.4
$ws.Range("L122").Value = 2085
$ws.Range("M122").Value = -3075.4
$ws.Range("N122").Value = -6985
$ws.Range("H137").Value = 2184.9092
$ws.Range("I137").Value = 1755.5
$ws.Range("J137").Value = 3330
$ws.Range("K137").Value = 5266.5
$ws.Range("L137").Value = 9990
$ws.Range("M137").Value = -2716.5
$ws.Range("N137").Value = -15090

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1099.8334
$ws.Range("I2").Value = 1099.8334
$ws.Range("K2").Value = 1099.8334
$ws.Range("M2").Value = -986.8334
$ws.Range("H32").Value = 9613.405000000001
$ws.Range("I32").Value = 9186.027
$ws.Range("J32").Value = 24999
$ws.Range("K32").Value = 9186.027
$ws.Range("L32").Value = 24999
$ws.Range("M32").Value = -8899.027
$ws.Range("N32").Value = -25573
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H74").Value = 4083.3333
$ws.Range("I74").Value = 4125
$ws.Range("K74").Value = 4125
$ws.Range("M74").Value = -3251
$ws.Range("H77").Value = 4083.3333
$ws.Range("I77").Value = 4125
$ws.Range("K77").Value = 20625
$ws.Range("M77").Value = -16257
$ws.Range("H116").Value = 1099.8334
$ws.Range("I116").Value = 1099.8334
$ws.Range("K116").Value = 1099.8334
$ws.Range("M116").Value = 1194.1666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1099.8334
$ws.Range("I3").Value = 1099.8334
$ws.Range("K3").Value = 1099.8334
$ws.Range("M3").Value = -985.8334
$ws.Range("H59").Value = 98998
$ws.Range("J59").Value = 98998
$ws.Range("L59").Value = 98998
$ws.Range("N59").Value = -100692
$ws.Range("H86").Value = 1900
$ws.Range("J86").Value = 1900
$ws.Range("L86").Value = 1900
$ws.Range("N86").Value = -4146
$ws.Range("H89").Value = 1900
$ws.Range("J89").Value = 1900
$ws.Range("L89").Value = 9500
$ws.Range("N89").Value = -20732
$ws.Range("H105").Value = 4861.5386
$ws.Range("I105").Value = 4150
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 4150
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -2403
$ws.Range("N105").Value = -9494
$ws.Range("H107").Value = 1457.8
$ws.Range("I107").Value = 1457.8
$ws.Range("K107").Value = 1457.8
$ws.Range("M107").Value = 462.2

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 862.4
$ws.Range("I22").Value = 983.2857
$ws.Range("J22").Value = 756.625
$ws.Range("K22").Value = 983.2857
$ws.Range("L22").Value = 756.625
$ws.Range("M22").Value = -633.2857
$ws.Range("N22").Value = -1456.625
$ws.Range("H31").Value = 4035.6667
$ws.Range("I31").Value = 3943.2
$ws.Range("K31").Value = 3943.2
$ws.Range("M31").Value = -3648.2
$ws.Range("H34").Value = 4035.6667
$ws.Range("I34").Value = 3943.2
$ws.Range("K34").Value = 3943.2
$ws.Range("M34").Value = -3741.2
$ws.Range("H60").Value = 6999.8335
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H132").Value = 5849.75
$ws.Range("I132").Value = 4542.5713
$ws.Range("K132").Value = 13627.7139
$ws.Range("M132").Value = -11097.7139
$ws.Range("H134").Value = 3876.8462
$ws.Range("I134").Value = 3458.7368
$ws.Range("K134").Value = 10376.2104
$ws.Range("M134").Value = -7841.2104

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 3000
$ws.Range("M111").Value = 67
$ws.Range("H113").Value = 975
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 975
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2925
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7265
$ws.Range("H128").Value = 499999
$ws.Range("I128").Value = 499999
$ws.Range("K128").Value = 1499997
$ws.Range("M128").Value = -1495017

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 275.125
$ws.Range("I2").Value = 21.3
$ws.Range("K2").Value = 21.3
$ws.Range("M2").Value = 91.7
$ws.Range("H52").Value = 225
$ws.Range("I52").Value = 225
$ws.Range("K52").Value = 225
$ws.Range("M52").Value = 34
$ws.Range("H102").Value = 2099.5
$ws.Range("I102").Value = 1200
$ws.Range("K102").Value = 1200
$ws.Range("M102").Value = 422
$ws.Range("H122").Value = 7414.7
$ws.Range("I122").Value = 8018.625
$ws.Range("K122").Value = 24055.875
$ws.Range("M122").Value = -21605.875
$ws.Range("H132").Value = 2865.45
$ws.Range("I132").Value = 2458.375
$ws.Range("K132").Value = 7375.125
$ws.Range("M132").Value = -4845.125

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1311.6
$ws.Range("I16").Value = 1309.5555
$ws.Range("J16").Value = 1330
$ws.Range("K16").Value = 1309.5555
$ws.Range("L16").Value = 1330
$ws.Range("M16").Value = -1139.5555
$ws.Range("N16").Value = -1670
$ws.Range("H61").Value = 1019
$ws.Range("I61").Value = 523.75
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 523.75
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -321.75
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 1019
$ws.Range("I113").Value = 523.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 523.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1646.25
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 4734.8335
$ws.Range("I132").Value = 2102.25
$ws.Range("K132").Value = 6306.75
$ws.Range("M132").Value = -3776.75
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
$ws.Range("H138").Value = 82857.14
$ws.Range("I138").Value = 80000
$ws.Range("J138").Value = 100000
$ws.Range("K138").Value = 80000
$ws.Range("L138").Value = 100000
$ws.Range("M138").Value = -74860
$ws.Range("N138").Value = -110280

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 24999
$ws.Range("J63").Value = 24999
$ws.Range("L63").Value = 24999
$ws.Range("N63").Value = -26247
$ws.Range("H66").Value = 24999
$ws.Range("J66").Value = 24999
$ws.Range("L66").Value = 74997
$ws.Range("N66").Value = -81237
$ws.Range("H113").Value = 977.7857
$ws.Range("I113").Value = 688.9
$ws.Range("K113").Value = 2066.7
$ws.Range("M113").Value = 103.3000000000002
$ws.Range("H132").Value = 2875.6667
$ws.Range("I132").Value = 1802.5555
$ws.Range("K132").Value = 5407.666499999999
$ws.Range("M132").Value = -2877.666499999999
